$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the stored credential values
$ws.Range("A2").Value = "venukollapudi@gmail.com"
$ws.Range("B2").Value = "Venu@12345"

# Add a mailto hyperlink on B2, matching the style used for A2
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:venukollapudi@gmail.com", "", "", "Venu@12345")

# Move / restore the active selection to A2
$ws.Range("A2").Select()
